$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: wrap text + vertical-center the data block A2:F5
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$rngSummary = $wsSummary.Range("A2:F5")
$rngSummary.WrapText = $true
$rngSummary.VerticalAlignment = -4108
$wsSummary.Range("A8").Select() | Out-Null

# ---------------------------------------------------------------------------
# Repayment schedule sheet: wrap text + vertical-center the data block A2:P8
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$rngSchedule = $wsSchedule.Range("A2:P8")
$rngSchedule.WrapText = $true
$rngSchedule.VerticalAlignment = -4108
$wsSchedule.Range("E3").Font.Italic = $true
$wsSchedule.Rows("9:9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Transactions sheet: corrected overdue test data + wrap/center formatting
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 46
$wsTransactions.Range("A3").Value = 44
$wsTransactions.Range("J2").Value = 4165.74

$rngRow2 = $wsTransactions.Range("A2:L2")
$rngRow2.WrapText = $true
$rngRow2.VerticalAlignment = -4108

$rngRow3 = $wsTransactions.Range("A3:J3")
$rngRow3.WrapText = $true
$rngRow3.VerticalAlignment = -4108

$wsTransactions.Range("K2:L2").Font.Italic = $true

# K3/L3 are touched but left with no special formatting (same as a fresh,
# un-formatted cell) - matches the source workbook which has bare <c/> nodes
$wsTransactions.Range("K3:L3").Style = "Normal"

$wsTransactions.Rows("2:5").Select() | Out-Null
$wsTransactions.Activate() | Out-Null
